
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix an existing data value: "Feedback Ticket Management" (row 13) Latest Version 16.0.0.1 -> 16.0.1.0
$ws.Range("D13").Value = "16.0.1.0"

# Insert a new row at 15 for the "grading_template" module, shifting rows 15.. down by one.
$ws.Rows.Item(15).Insert()

$ws.Range("A15").Value = "grading_template"
$ws.Range("B15").Value = "My Company"
$ws.Range("C15").Value = "https://www.yourcompany.com"
$ws.Range("D15").Value = "16.0.0.1"
$ws.Range("E15").Value = "Installed"

$ws.Hyperlinks.Add($ws.Range("C15"), "https://www.yourcompany.com")
